$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Foglio1")

# Add new column header "VQ_best" in T1
$ws.Range("T1").Value = "VQ_best"

# Add new value in T2
$ws.Range("T2").Value = 1

# Update selection to match the saved view state
$ws.Activate()
$ws.Range("T2").Select()
